$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.000.11"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "1.643.15"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "214.74"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.2565"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "0.06361"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").Value = "0.07766"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.281"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.634.38"
$ws.Range("E13").Value = "  -1.05%  "

$ws.Range("D14").Value = "0.5437"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "64.29"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₅7726"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").Value = "26.025.63"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "196.99"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").Value = "4.415"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").Value = "9.923"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").Value = "6.029"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "1.865"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").Value = "141.48"
$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("D26").Value = "0.1190"

$ws.Range("D27").Value = "6.820"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").Value = "1.234"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("D30").Value = "0.04852"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").Value = "3.251"
$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("D32").Value = "3.162"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").Value = "2.366"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").Value = "0.8960"
$ws.Range("E35").Value = "  +1.01%  "

$ws.Range("D36").Value = "2.579"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").Value = "1.141.23"
$ws.Range("E37").Value = "  +1.01%  "

$ws.Range("D38").Value = "0.5443"

$ws.Range("D39").Value = "0.01562"
$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "2.522"
$ws.Range("E41").Value = "  -1.80%  "

$ws.Range("E42").Value = "  +9.22%  "

$ws.Range("D43").Value = "0.8097"
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").Value = "99.23"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").Value = "5.397"
$ws.Range("E45").Value = "  -4.97%  "

$ws.Range("D46").Value = "1.779.28"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").Value = "0.4530"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "0.9998"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "54.83"
$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").Value = "0.05057"
$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.29%  "

